# "work on checkout page and other pages created"
#
# Checkout page (slide 20, "Scan barcode" / checkout UI) got rearranged:
#   - "Or input manually" textbox moved further up the page
#   - "Add to basket" textbox moved further up the page
#   - the now-redundant "Details of item after barcode has been scanned"
#     helper textbox was deleted
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)

$shapes = $s.Shapes

# walk backwards so deleting a shape doesn't shift the index of shapes
# still to be visited
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $shp = $shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $text = $shp.TextFrame.TextRange.Text

    if ($text -eq "Or input manually") {
        # 2163516,5095522 -> 2212606,4181122 EMU
        $shp.Left = 174.22094488188975
        $shp.Top = 329.2222137454094
    }
    elseif ($text -eq "Add to basket") {
        # 2354250,4296279 -> 2402566,3390206 EMU
        $shp.Left = 189.17842865085038
        $shp.Top = 266.94535827670865
    }
    elseif ($text -eq "Details of item after barcode has been scanned") {
        $shp.Delete()
    }
}
